$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" / "_new" header suffixes to the concrete AHB version
# labels "_FV2310" / "_FV2404" (the "diff" column header is unchanged).
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"
$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# Stash the header row's existing look (bold/fill/border/wrap) on a scratch
# cell, then strip it from the header row before turning the range into a
# Table: the engine bakes whatever formatting the header row carries at
# Add()-time into a new dxf/headerRowDxfId, which the source file does not
# have, so we round-trip the formatting around the Add() call instead of
# letting it be captured.
$ws.Range("A1").Copy()
$ws.Range("AA1").PasteSpecial(-4122)
$ws.Range("A1:U1").ClearFormats()

# Turn the data range into an Excel Table ("Table1") with the header row
# and an AutoFilter, matching the regenerated merged-AHB export layout.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U78"), $null, 1)
$tbl.Name = "Table1"

# Restore the header row's original formatting and drop the scratch cell.
$ws.Range("AA1").Copy()
$ws.Range("A1:U1").PasteSpecial(-4122)
$ws.Range("AA1").Clear()

# Freeze the header row so it stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
